# Dungeon.xlsx : add a "Jobs" column (int[] / job id list) so players can
# select which job to play when entering a dungeon.
# The new column is inserted right after "CardDeck" (old column H),
# pushing Str/Agi/Intl/Perc/Endu/QuestDungeon/QuestDungeonRate/BgImage one
# column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new worksheet column at H; this shifts the existing H:O data
# (Str..BgImage) right to I:P, carrying values/styles along with it.
$ws.Range("H1").EntireColumn.Insert()

# Grow the table definition so it covers the new column.
$lo.Resize($ws.Range("A3:P7"))

# --- Populate the new "Jobs" column -----------------------------------
# Row 1: Chinese display name; Row 2: data type; Row 3: field/header name;
# Rows 4-7: the actual job-id list value for every dungeon.
$ws.Range("H1").Value = "职业列表"
$ws.Range("H2").Value = "int[]"
$ws.Range("H3").Value = "Jobs"

$ws.Range("H4").Value = "11000001;11000002;11000003"
$ws.Range("H5").Value = "11000001;11000002;11000003"
$ws.Range("H6").Value = "11000001;11000002;11000003"
$ws.Range("H7").Value = "11000001;11000002;11000003"

# Match styling used by the neighbouring header/data cells.
$ws.Range("H1").Style = $ws.Range("G1").Style
$ws.Range("H2").Style = $ws.Range("G2").Style
$ws.Range("H3").Style = $ws.Range("G3").Style
$ws.Range("H4:H7").Style = $ws.Range("G4:G7").Style

# Match the selection left behind by the original author's edit.
$ws.Range("H4:H7").Select()
